$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '96.641.39'
$ws.Range('E2').Value = '  +4.79%  '
$ws.Range('D3').Value = '3.101.02'
$ws.Range('E3').Value = '  +0.13%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = '''236.23'
$ws.Range('E5').Value = '  +1.20%  '
$ws.Range('D6').Value = '''601.85'
$ws.Range('E6').Value = '  -1.55%  '
$ws.Range('D7').Value = '''1.09'
$ws.Range('E7').Value = '  +1.14%  '
$ws.Range('E8').Value = '  -1.95%  '
$ws.Range('D9').Value = '''1.00'
$ws.Range('E9').Value = '  +0.11%  '
$ws.Range('D10').Value = '3.096.22'
$ws.Range('E10').Value = '  +0.08%  '
$ws.Range('E11').Value = '  +0.53%  '
$ws.Range('D12').Value = '''0.196'
$ws.Range('E12').Value = '  -0.38%  '
$ws.Range('D13').Value = '95.541.35'
$ws.Range('E13').Value = '  +3.88%  '
$ws.Range('E14').Value = '  -3.02%  '
$ws.Range('D15').Value = '''33.24'
$ws.Range('E15').Value = '  -0.94%  '
$ws.Range('D16').Value = '''5.28'
$ws.Range('E16').Value = '  -1.77%  '
$ws.Range('D17').Value = '3.669.82'
$ws.Range('E17').Value = '  -0.08%  '
$ws.Range('D18').Value = '3.074.92'
$ws.Range('E18').Value = '  -0.39%  '
$ws.Range('E19').Value = '  -9.02%  '
$ws.Range('D20').Value = '''14.23'
$ws.Range('E20').Value = '  -0.79%  '
$ws.Range('D21').Value = '''464.39'
$ws.Range('E21').Value = '  +7.07%  '
$ws.Range('D22').Value = '''5.57'
$ws.Range('E22').Value = '  -2.84%  '
$ws.Range('B23').Value = 'Uniswap'
$ws.Range('C23').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D23').Value = '''8.65'
$ws.Range('E23').Value = '  -4.28%  '
$ws.Range('B24').Value = 'PEPE'
$ws.Range('C24').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D24').Value = '''0.0000185'
$ws.Range('E24').Value = '  -6.53%  '
$ws.Range('D25').Value = '''5.44'
$ws.Range('E25').Value = '  -1.98%  '
$ws.Range('D26').Value = '''85.05'
$ws.Range('E26').Value = '  +0.04%  '
$ws.Range('D27').Value = '''11.39'
$ws.Range('E27').Value = '  +0.69%  '
$ws.Range('D28').Value = '3.252.14'
$ws.Range('E29').Value = '  -0.03%  '
$ws.Range('D30').Value = '''0.175'
$ws.Range('E30').Value = '  -1.30%  '
$ws.Range('D31').Value = '''0.235'
$ws.Range('E31').Value = '  +1.33%  '
$ws.Range('D32').Value = '''0.124'
$ws.Range('E32').Value = '  +0.18%  '
$ws.Range('D33').Value = '''8.86'
$ws.Range('E33').Value = '  -2.41%  '
$ws.Range('D34').Value = '''25.86'
$ws.Range('E34').Value = '  +1.89%  '
$ws.Range('D35').Value = '''0.823'
$ws.Range('E35').Value = '  -20.58%  '
$ws.Range('B36').Value = 'Kaspa'
$ws.Range('C36').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D36').Value = '''0.148'
$ws.Range('E36').Value = '  -3.72%  '
$ws.Range('B37').Value = 'RenderToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D37').Value = '''7.20'
$ws.Range('E37').Value = '  -8.76%  '
$ws.Range('B38').Value = 'Bittensor'
$ws.Range('C38').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D38').Value = '''480.66'
$ws.Range('E38').Value = '  +3.54%  '
$ws.Range('B39').Value = 'WhiteBITCoin'
$ws.Range('C39').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D39').Value = '''24.16'
$ws.Range('E39').Value = '  +1.34%  '
$ws.Range('E40').Value = '  -3.18%  '
$ws.Range('D41').Value = '''0.430'
$ws.Range('E41').Value = '  -1.06%  '
$ws.Range('D42').Value = '''3.62'
$ws.Range('E42').Value = '  -6.57%  '
$ws.Range('E43').Value = '  -4.58%  '
$ws.Range('E44').Value = '  -0.08%  '
$ws.Range('D45').Value = '''3.08'
$ws.Range('E45').Value = '  -6.12%  '
$ws.Range('D46').Value = '''161.75'
$ws.Range('E46').Value = '  +1.88%  '
$ws.Range('B47').Value = 'Stacks'
$ws.Range('C47').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D47').Value = '''1.85'
$ws.Range('E47').Value = '  +1.84%  '
$ws.Range('B48').Value = 'ARBITRUM'
$ws.Range('C48').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D48').Value = '''0.676'
$ws.Range('E48').Value = '  +0.13%  '
$ws.Range('B49').Value = 'OKB'
$ws.Range('C49').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D49').Value = '''43.79'
$ws.Range('E49').Value = '  +0.04%  '
$ws.Range('B50').Value = 'FirstDigitalUSD'
$ws.Range('C50').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D50').Value = '''0.999'
$ws.Range('E50').Value = '  +0.14%  '
$ws.Range('B51').Value = 'FLOKI'
$ws.Range('C51').Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range('D51').Value = '''0.000264'
$ws.Range('E51').Value = '  +9.11%  '
